# Apply the data updates described by the diff.
# Sheet 1 = "展览" (Exhibition)
# Sheet 2 = "演出" (Performance)
# Sheet 3 = "本地生活" (Local life) -- unchanged
# Sheet 4 = "全部类型" (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7715
$ws1.Range("F4").Value  = 7887
$ws1.Range("F7").Value  = 33
$ws1.Range("F8").Value  = 6678
$ws1.Range("F9").Value  = 6678
$ws1.Range("F10").Value = 3388
$ws1.Range("F21").Value = 317
$ws1.Range("F23").Value = 331
$ws1.Range("F24").Value = 3862
$ws1.Range("F32").Value = 2764
$ws1.Range("F33").Value = 1846
$ws1.Range("F37").Value = 3699
$ws1.Range("F38").Value = 324
$ws1.Range("F41").Value = 921
$ws1.Range("F47").Value = 553
$ws1.Range("F48").Value = 643
$ws1.Range("F49").Value = 2
$ws1.Range("G49").Value = 70

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value  = 680
$ws2.Range("F6").Value  = 412
$ws2.Range("F17").Value = 37

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G4").Value  = 680
$ws4.Range("F7").Value  = 7715
$ws4.Range("F9").Value  = 7887
$ws4.Range("F11").Value = 33
$ws4.Range("F12").Value = 6678
$ws4.Range("F13").Value = 3388
$ws4.Range("F24").Value = 317
$ws4.Range("F25").Value = 331
$ws4.Range("F26").Value = 3862
$ws4.Range("F35").Value = 2764
$ws4.Range("F36").Value = 1846
$ws4.Range("F41").Value = 324
$ws4.Range("F45").Value = 921
$ws4.Range("F47").Value = 37
$ws4.Range("F49").Value = 553
$ws4.Range("F50").Value = 643
